# "add 2d act camera"
# The villageScene row (row 2, ID=1) gets updated 2D-act camera offset
# position/rotation values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"

[void]$ws.Range("K2").Select()
